$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "81.161.86"
$ws.Range("E2").Value = "  +2.73%  "
$ws.Range("D3").Value = "3.142.55"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "618.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.284"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +24.25%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.582"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("D10").Value = "3.137.59"
$ws.Range("E10").Value = "  -1.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.584"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.77%  "
$ws.Range("E12").Value = "  +11.37%  "
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.27"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.28%  "
$ws.Range("D15").Value = "3.713.87"
$ws.Range("E15").Value = "  -1.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "31.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("D17").Value = "81.070.32"
$ws.Range("E17").Value = "  +2.74%  "
$ws.Range("D18").Value = "3.131.89"
$ws.Range("E18").Value = "  -1.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "430.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.91"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.60%  "
$ws.Range("D26").Value = "3.301.80"
$ws.Range("E26").Value = "  -1.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "75.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.67%  "
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("E30").Value = "  +5.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "581.48"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +11.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.154"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +13.93%  "
$ws.Range("E36").Value = "  -0.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.137"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +11.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "22.64"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.407"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.22%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +10.60%  "
$ws.Range("B42").Value = "WhiteBITCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "20.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.51%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.04"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +21.40%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +12.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "158.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.55%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "185.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "45.02"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.31"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.767"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.03%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "25.90"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.40%  "
